{"js": "// Word adds a new acronym row (\"<NombreResultado>\" / \"Nombre del resultado\n// del sprint\") to the \"Sigla / Significado\" table, placed immediately\n// before the existing \"<ITERXX>\" row.\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// The acronym (\"Sigla / Significado\") table is the second table in the\n// document (2 columns). Find it defensively by column count in case the\n// document layout shifts.\nlet siglaTable = null;\nfor (const t of tables.items) {\n  t.load(\"headerRowCount\");\n}\nawait context.sync();\n\nfor (const t of tables.items) {\n  const rows = t.rows;\n  rows.load(\"items\");\n}\nawait context.sync();\n\nfor (const t of tables.items) {\n  const firstRow = t.rows.items[0];\n  firstRow.load(\"cellCount\");\n}\nawait context.sync();\n\nfor (const t of tables.items) {\n  if (t.rows.items.length > 0 && t.rows.items[0].cellCount === 2) {\n    siglaTable = t;\n    break;\n  }\n}\nif (!siglaTable) {\n  siglaTable = tables.items[tables.items.length - 1];\n}\n\n// Locate the \"<ITERXX>\" row within that table so the new row can be\n// inserted right before it.\nconst rows = siglaTable.rows;\nrows.load(\"items\");\nawait context.sync();\nrows.items.forEach((r) => r.load(\"values\"));\nawait context.sync();\n\nlet iterRowIndex = -1;\nfor (let i = 0; i < rows.items.length; i++) {\n  const values = rows.items[i].values;\n  if (values && values[0] && values[0][0] && values[0][0].indexOf(\"<ITERXX>\") !== -1) {\n    iterRowIndex = i;\n    break;\n  }\n}\n\nif (iterRowIndex === -1) {\n  throw new Error('Could not locate the \"<ITERXX>\" row in the acronym table.');\n}\n\nconst iterRow = rows.items[iterRowIndex];\niterRow.insertRows(\"Before\", 1, [[\"<NombreResultado>\", \"Nombre del resultado del sprint\"]]);\nawait context.sync();\n", "ps1": "# Word adds a new acronym row (\"<NombreResultado>\" / \"Nombre del resultado\n# del sprint\") to the \"Sigla / Significado\" table, placed immediately\n# before the existing \"<ITERXX>\" row. The table's first-column cell width\n# is also nudged from 4513 dxa (225.65pt) to 4514 dxa (225.7pt).\n\n$d = $word.ActiveDocument\n\n# Find the acronym (\"Sigla / Significado\") table - the 2-column table that\n# contains the \"<ITERXX>\" placeholder.\n$tbl = $null\nfor ($t = 1; $t -le $d.Tables.Count; $t++) {\n    $candidate = $d.Tables.Item($t)\n    if ($candidate.Columns.Count -eq 2) {\n        $tbl = $candidate\n        break\n    }\n}\nif ($tbl -eq $null) {\n    $tbl = $d.Tables.Item($d.Tables.Count)\n}\n\n# Locate the row whose first cell holds \"<ITERXX>\" so the new row can be\n# inserted right before it.\n$iterRowIndex = -1\nfor ($r = 1; $r -le $tbl.Rows.Count; $r++) {\n    $cellText = $tbl.Cell($r, 1).Range.Text\n    if ($cellText -like \"*<ITERXX>*\") {\n        $iterRowIndex = $r\n        break\n    }\n}\nif ($iterRowIndex -eq -1) {\n    throw \"Could not locate the '<ITERXX>' row in the acronym table.\"\n}\n\n$iterRow = $tbl.Rows.Item($iterRowIndex)\n$newRow = $tbl.Rows.Add($iterRow)\n$newRow.Cells.Item(1).Range.Text = \"<NombreResultado>\"\n$newRow.Cells.Item(2).Range.Text = \"Nombre del resultado del sprint\"\n\n# Fix the first column's stored cell width (4513 -> 4514 dxa). Setting the\n# width on any cell in the column re-applies it to the whole column.\n$tbl.Cell(1, 1).Width = 225.7\n"}
